$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.4636016765006184
$ws.Range("E2").Value = 0.5293505787849426
$ws.Range("F2").Value = 0.5399087732845833
$ws.Range("G2").Value = 0.45448322019344
$ws.Range("H2").Value = 24347000000
$ws.Range("I2").Value = "NVDA"

$ws.Range("D3").Value = 0.5390969049603723
$ws.Range("E3").Value = 0.5330748558044434
$ws.Range("F3").Value = 0.5511411181343596
$ws.Range("G3").Value = 0.4909202819886819
$ws.Range("H3").Value = 24347000000
$ws.Range("I3").Value = "NVDA"

$ws.Range("D4").Value = 0.4866588576917522
$ws.Range("E4").Value = 0.5440124869346619
$ws.Range("F4").Value = 0.5795862766906059
$ws.Range("G4").Value = 0.4803669213712621
$ws.Range("H4").Value = 24347000000
$ws.Range("I4").Value = "NVDA"

$ws.Range("D5").Value = 0.6916250665320051
$ws.Range("E5").Value = 0.7708484530448914
$ws.Range("F5").Value = 0.7761948350357134
$ws.Range("G5").Value = 0.671454707027059
$ws.Range("H5").Value = 24347000000
$ws.Range("I5").Value = "NVDA"

$ws.Range("D6").Value = 0.7137665718164842
$ws.Range("E6").Value = 0.7649937868118286
$ws.Range("F6").Value = 0.7835331714896975
$ws.Range("G6").Value = 0.6037498721048875
$ws.Range("H6").Value = 24347000000
$ws.Range("I6").Value = "NVDA"

$ws.Range("D7").Value = 0.8781987265865948
$ws.Range("E7").Value = 1.143837332725525
$ws.Range("F7").Value = 1.146285670760405
$ws.Range("G7").Value = 0.8422090213734176
$ws.Range("H7").Value = 24347000000
$ws.Range("I7").Value = "NVDA"

$ws.Range("D8").Value = 1.407489974210975
$ws.Range("E8").Value = 1.505675435066223
$ws.Range("F8").Value = 1.558695654155838
$ws.Range("G8").Value = 1.36232461539904
$ws.Range("H8").Value = 24347000000
$ws.Range("I8").Value = "NVDA"

$ws.Range("D9").Value = 1.756334030533547
$ws.Range("E9").Value = 2.267348051071167
$ws.Range("F9").Value = 2.342352353151048
$ws.Range("G9").Value = 1.637310506455994
$ws.Range("H9").Value = 24347000000
$ws.Range("I9").Value = "NVDA"

$ws.Range("D10").Value = 2.72341255007313
$ws.Range("E10").Value = 2.499293804168701
$ws.Range("F10").Value = 2.978070694070083
$ws.Range("G10").Value = 2.356941417088456
$ws.Range("H10").Value = 24347000000
$ws.Range("I10").Value = "NVDA"

$ws.Range("D11").Value = 2.582837923141818
$ws.Range("E11").Value = 3.559601545333862
$ws.Range("F11").Value = 3.624949117864806
$ws.Range("G11").Value = 2.522915286563916
$ws.Range("H11").Value = 24347000000
$ws.Range("I11").Value = "NVDA"

$ws.Range("D12").Value = 4.00225989959374
$ws.Range("E12").Value = 4.182711124420166
$ws.Range("F12").Value = 4.309100668030247
$ws.Range("G12").Value = 3.774659845002095
$ws.Range("H12").Value = 24347000000
$ws.Range("I12").Value = "NVDA"

$ws.Range("D13").Value = 5.172457752856373
$ws.Range("E13").Value = 4.958987236022949
$ws.Range("F13").Value = 5.402729379806347
$ws.Range("G13").Value = 4.724762739341751
$ws.Range("H13").Value = 24347000000
$ws.Range("I13").Value = "NVDA"

$ws.Range("D14").Value = 5.897262289691401
$ws.Range("E14").Value = 5.983303546905518
$ws.Range("F14").Value = 6.22980569137383
$ws.Range("G14").Value = 5.043776290344111
$ws.Range("H14").Value = 24347000000
$ws.Range("I14").Value = "NVDA"

$ws.Range("D15").Value = 5.555807205017004
$ws.Range("E15").Value = 6.239119052886963
$ws.Range("F15").Value = 6.444706271171777
$ws.Range("G15").Value = 5.497173632023759
$ws.Range("H15").Value = 24347000000
$ws.Range("I15").Value = "NVDA"

$ws.Range("D16").Value = 6.092963007280932
$ws.Range("E16").Value = 6.948250293731689
$ws.Range("F16").Value = 6.973995578263955
$ws.Range("G16").Value = 5.909527736947218
$ws.Range("H16").Value = 24347000000
$ws.Range("I16").Value = "NVDA"

$ws.Range("D17").Value = 5.258371480907623
$ws.Range("E17").Value = 4.047930717468262
$ws.Range("F17").Value = 5.498627103393312
$ws.Range("G17").Value = 3.301900774778706
$ws.Range("H17").Value = 24347000000
$ws.Range("I17").Value = "NVDA"

$ws.Range("D18").Value = 3.582643338630164
$ws.Range("E18").Value = 3.824626684188842
$ws.Range("F18").Value = 4.097849788591541
$ws.Range("G18").Value = 3.535040189327395
$ws.Range("H18").Value = 24347000000
$ws.Range("I18").Value = "NVDA"

$ws.Range("D19").Value = 4.544595380371072
$ws.Range("E19").Value = 3.36197280883789
$ws.Range("F19").Value = 4.588276380167989
$ws.Range("G19").Value = 3.360235490966609
$ws.Range("H19").Value = 24347000000
$ws.Range("I19").Value = "NVDA"

$ws.Range("D20").Value = 4.202666232865935
$ws.Range("E20").Value = 4.16216516494751
$ws.Range("F20").Value = 4.30975788284301
$ws.Range("G20").Value = 3.662238416308913
$ws.Range("H20").Value = 24347000000
$ws.Range("I20").Value = "NVDA"

$ws.Range("D21").Value = 4.964422297274594
$ws.Range("E21").Value = 5.390726089477539
$ws.Range("F21").Value = 5.506877897668168
$ws.Range("G21").Value = 4.939301764120212
$ws.Range("H21").Value = 24347000000
$ws.Range("I21").Value = "NVDA"

$ws.Range("D22").Value = 5.866621848020344
$ws.Range("E22").Value = 6.722098350524902
$ws.Range("F22").Value = 7.873270450584467
$ws.Range("G22").Value = 5.860399391076736
$ws.Range("H22").Value = 24347000000
$ws.Range("I22").Value = "NVDA"

$ws.Range("D23").Value = 7.08176428071292
$ws.Range("E23").Value = 8.84180736541748
$ws.Range("F23").Value = 9.14689537722265
$ws.Range("G23").Value = 6.994347431035929
$ws.Range("H23").Value = 24347000000
$ws.Range("I23").Value = "NVDA"

$ws.Range("D24").Value = 10.69663803758115
$ws.Range("E24").Value = 13.32981014251709
$ws.Range("F24").Value = 13.52963971933032
$ws.Range("G24").Value = 10.67944561664404
$ws.Range("H24").Value = 24347000000
$ws.Range("I24").Value = "NVDA"

$ws.Range("D25").Value = 12.61922874996378
$ws.Range("E25").Value = 13.36071491241455
$ws.Range("F25").Value = 14.64678901241372
$ws.Range("G25").Value = 12.35727870511663
$ws.Range("H25").Value = 24347000000
$ws.Range("I25").Value = "NVDA"

$ws.Range("D26").Value = 13.01737104254144
$ws.Range("E26").Value = 13.67680358886719
$ws.Range("F26").Value = 15.330246261919
$ws.Range("G26").Value = 12.86728442287
$ws.Range("H26").Value = 24347000000
$ws.Range("I26").Value = "NVDA"

$ws.Range("D27").Value = 15.08863484284466
$ws.Range("E27").Value = 16.20544242858887
$ws.Range("F27").Value = 16.23836403019814
$ws.Range("G27").Value = 13.42664005449409
$ws.Range("H27").Value = 24347000000
$ws.Range("I27").Value = "NVDA"

$ws.Range("D28").Value = 19.65713514402224
$ws.Range("E28").Value = 22.3362922668457
$ws.Range("F28").Value = 22.99285955249893
$ws.Range("G28").Value = 18.72117452752336
$ws.Range("H28").Value = 24347000000
$ws.Range("I28").Value = "NVDA"

$ws.Range("D29").Value = 25.59771005461229
$ws.Range("E29").Value = 32.61065673828125
$ws.Range("F29").Value = 34.5777160144105
$ws.Range("G29").Value = 25.17655304160672
$ws.Range("H29").Value = 24347000000
$ws.Range("I29").Value = "NVDA"

$ws.Range("D30").Value = 25.05685828421552
$ws.Range("E30").Value = 24.33920860290528
$ws.Range("F30").Value = 26.87443784756249
$ws.Range("G30").Value = 20.8507708873018
$ws.Range("H30").Value = 24347000000
$ws.Range("I30").Value = "NVDA"

$ws.Range("D31").Value = 18.50933892148411
$ws.Range("E31").Value = 18.64011573791504
$ws.Range("F31").Value = 20.36516366321006
$ws.Range("G31").Value = 15.5404174802675
$ws.Range("H31").Value = 24347000000
$ws.Range("I31").Value = "NVDA"

$ws.Range("D32").Value = 18.15478572494654
$ws.Range("E32").Value = 15.07140827178955
$ws.Range("F32").Value = 19.24515225937969
$ws.Range("G32").Value = 14.93661010133938
$ws.Range("H32").Value = 24347000000
$ws.Range("I32").Value = "NVDA"

$ws.Range("D33").Value = 13.79442516241902
$ws.Range("E33").Value = 16.90269088745117
$ws.Range("F33").Value = 16.97759973622608
$ws.Range("G33").Value = 12.94045171299164
$ws.Range("H33").Value = 24347000000
$ws.Range("I33").Value = "NVDA"

$ws.Range("D34").Value = 19.67240398882402
$ws.Range("E34").Value = 23.1940746307373
$ws.Range("F34").Value = 23.86544088440668
$ws.Range("G34").Value = 19.59247961629646
$ws.Range("H34").Value = 24347000000
$ws.Range("I34").Value = "NVDA"

$ws.Range("D35").Value = 27.81842812694343
$ws.Range("E35").Value = 37.80468368530273
$ws.Range("F35").Value = 41.90550376267679
$ws.Range("G35").Value = 27.21889266011253
$ws.Range("H35").Value = 24347000000
$ws.Range("I35").Value = "NVDA"

$ws.Range("D36").Value = 46.42881177243578
$ws.Range("E36").Value = 49.32186889648438
$ws.Range("F36").Value = 50.2322566659547
$ws.Range("G36").Value = 40.28394116435364
$ws.Range("H36").Value = 24347000000
$ws.Range("I36").Value = "NVDA"

$ws.Range("D37").Value = 40.85991626068604
$ws.Range("E37").Value = 46.74245071411133
$ws.Range("F37").Value = 50.51822517131968
$ws.Range("G37").Value = 40.84492570638794
$ws.Range("H37").Value = 24347000000
$ws.Range("I37").Value = "NVDA"

$ws.Range("D38").Value = 62.06887471690946
$ws.Range("E38").Value = 79.07234954833984
$ws.Range("F38").Value = 82.35270233203946
$ws.Range("G38").Value = 61.61910330161587
$ws.Range("H38").Value = 24347000000
$ws.Range("I38").Value = "NVDA"

$ws.Range("D39").Value = 85.03835422698414
$ws.Range("E39").Value = 109.5831985473633
$ws.Range("F39").Value = 115.7663853476814
$ws.Range("G39").Value = 81.21808428354018
$ws.Range("H39").Value = 24347000000
$ws.Range("I39").Value = "NVDA"

$ws.Range("D40").Value = 117.486257609687
$ws.Range("E40").Value = 119.3255767822266
$ws.Range("F40").Value = 131.2111434412157
$ws.Range("G40").Value = 90.6562503200763
$ws.Range("H40").Value = 24347000000
$ws.Range("I40").Value = "NVDA"

$ws.Range("D41").Value = 134.6613874226879
$ws.Range("E41").Value = 138.2103729248047
$ws.Range("F41").Value = 152.8461760008425
$ws.Range("G41").Value = 131.7622247614744
$ws.Range("H41").Value = 24347000000
$ws.Range("I41").Value = "NVDA"

$ws.Range("D42").Value = 114.725008510696
$ws.Range("E42").Value = 124.8927917480469
$ws.Range("F42").Value = 143.4087625347676
$ws.Range("G42").Value = 112.9853896023774
$ws.Range("H42").Value = 24347000000
$ws.Range("I42").Value = "NVDA"

$ws.Range("D43").Value = 113.065777256869
$ws.Range("E43").Value = 135.1130065917969
$ws.Range("F43").Value = 143.4719555798726
$ws.Range("G43").Value = 110.8060594111189
$ws.Range("H43").Value = 24347000000
$ws.Range("I43").Value = "NVDA"

$ws.Range("D44").Value = 174.0801747550306
$ws.Range("E44").Value = 174.170166015625
$ws.Range("F44").Value = 184.4695879752146
$ws.Range("G44").Value = 168.7904799128549
$ws.Range("H44").Value = 24347000000
$ws.Range("I44").Value = "NVDA"
